$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60.5342651418433
$ws.Range("L2").Value = 106.101642733938
$ws.Range("M2").Value = 25.3627907998722

$ws.Range("H4").Value = 4.93007252540714
$ws.Range("I4").Value = 8.20702036594837
$ws.Range("J4").Value = 6.39250428517991

$ws.Range("M13").Value = 0.46560922202107

$ws.Range("B19").Value = -11217.3674270217
$ws.Range("D19").Value = 0.746733701835474
$ws.Range("L19").Value = 106.101642733938
